$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New review / comment notes added off to the side of the plan
$ws.Range("P4").Value = "Gui for MPI"
$ws.Range("P5").Value = "Future cluster wrapper"
$ws.Range("P6").Value = "C++ framework optimization"
$ws.Range("P3").Value = "mpiexec configurations for non-local clusters"
$ws.Range("P2").Value = "Remains on 01/11/2020 in Parallel C++ messagin"
$ws.Range("U3").Value = "'- minimal requested to state its completed"
$ws.Range("U5").Value = "'- needed to run it on SCARF or future iDaaaS cluster"

# Header row: rename task placeholders from "W#" to "M#"
$ws.Range("A1").Value = "M1"
$ws.Range("B1").Value = "M2"
$ws.Range("C1").Value = "M2"
$ws.Range("D1").Value = "M4"
$ws.Range("E1").Value = "M5"
$ws.Range("F1").Value = "M6"
$ws.Range("G1").Value = "M7"
$ws.Range("H1").Value = "M8"
$ws.Range("I1").Value = "M9"
$ws.Range("J1").Value = "M10"
$ws.Range("K1").Value = "M11"
$ws.Range("L1").Value = "M12"

# Leave the selection where the author last clicked
$ws.Range("Q7").Select()
